$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "James Calderon, poor"
$ws.Range("B2").Value = "Niko Morris, poor"
$ws.Range("C2").Value = "Violet Hudson, poor"
$ws.Range("D2").Value = "Ava Lee, poor"
$ws.Range("E2").Value = "Caitlin Boyd, poor"
$ws.Range("F2").Value = "Stanley Hirst, poor"
$ws.Range("G2").Value = "Madison Taylor, good"
$ws.Range("H2").Value = "William Hunt, good"
$ws.Range("I2").Value = "Lexi Green, poor"
$ws.Range("J2").Value = "James Shilton, poor"

# Row 3 (A3 stays empty/unchanged)
$ws.Range("B3").Value = "Katrina Petersone, good"
$ws.Range("C3").Value = "Ruby Haigh, good"
$ws.Range("D3").Value = "Thomas Barrett, excellent"
$ws.Range("E3").Value = "Alex Sentance, excellent"
$ws.Range("F3").Value = "Nancy Enyoazu, good"
$ws.Range("G3").Value = "Benjamin Finn, good"
$ws.Range("H3").Value = "Brooke Layton, good"
$ws.Range("I3").Value = "Aarron Kelly, good"
$ws.Range("J3").Value = "Benedict Hobday, good"

# Row 4 (A4, I4, J4 stay empty/unchanged)
$ws.Range("B4").Value = "Esther Sido, excellent"
$ws.Range("C4").Value = "Sophie Rayner, excellent"
$ws.Range("D4").Value = "Benjamin Hillary, excellent"
$ws.Range("E4").Value = "James Eilbeck, excellent"
$ws.Range("F4").Value = "Spencer Rowe, excellent"
$ws.Range("G4").Value = "Samuel Dixon, excellent"
$ws.Range("H4").Value = "Matthew Homan, excellent"

$wb.Save()
